# "Fixed headers in test data"
# The "Connections" sheet header row used inconsistent / stale header names
# (OriginChain, OriginProcess, DestinationChain). Rename them to match the
# naming convention used elsewhere (Origin_Chain, Origin_Unit, Destination_Chain).

$wb = $excel.ActiveWorkbook

$connections = $wb.Worksheets.Item("Connections")

# Fix the header labels on the "Connections" sheet.
$connections.Range("A1").Value = "Origin_Chain"
$connections.Range("B1").Value = "Origin_Unit"
$connections.Range("F1").Value = "Destination_Chain"

# The workbook was left with the "Connections" tab active/selected (instead
# of "Chain List"), with the cursor sitting on C30.
$connections.Activate()
$connections.Range("C30").Select()

Write-Output "Headers fixed on 'Connections' sheet"
